# 4.8.xlsx update:
#  - rename two header cells (L1, M1) and make them bold
#  - widen columns L and M to fit the new, longer header text
#  - set the printed page setup (paper size / orientation)
#  - leave the last active selection on N8 (cosmetic, matches the saved file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Деньги (1)" / "Комиссионные (2)" headers to the new wording
# and make them bold (this introduces a new bold font + a new cell style
# that reuses the existing header fill/border/alignment).
$ws.Range("L1").Value() = "Деньги (x = 0)"
$ws.Range("M1").Value() = "Комиссионные (x = 1)"
$ws.Range("L1:M1").Font.Bold() = $true

# The longer bold labels need wider columns.
$ws.Range("L1").ColumnWidth() = 13.67
$ws.Range("M1").ColumnWidth() = 20.67

# Configure the print setup (A4, portrait).
$ws.PageSetup.PaperSize() = 9
$ws.PageSetup.Orientation() = 1

# Restore the last selected cell as saved in the workbook.
$ws.Range("N8").Select()
